{"js": "// Apply the benchmark-table numeric updates described by the commit:\n// a handful of single-value cells get new numbers, and three cells that\n// used to hold a full tab-separated row of stats are collapsed down to\n// just their first (leading) number.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Replace the entire text of a single-column table cell (row/col are\n// zero-based) while preserving the existing run formatting (font, size)\n// by reusing the paragraph's own Range instead of inserting a brand new\n// default-formatted run.\nasync function setCellText(rowIndex, text) {\n  const cell = table.getCell(rowIndex, 0);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  const range = paragraphs.items[0].getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\n// rowIndex (0-based) -> new text\nconst edits = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"881\"],\n  [6, \"0.05588\"],\n  [7, \"0.03023\"],\n  [9, \"0.11146\"],\n  [11, \"59.86878\"],\n  [43, \"75.98\"],\n  [44, \"59.87\"],\n  [45, \"249\"],\n];\n\nfor (const [rowIndex, text] of edits) {\n  await setCellText(rowIndex, text);\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-table numeric updates described by the commit:\n# a handful of single-value cells get new numbers, and three cells that\n# used to hold a full tab-separated row of stats are collapsed down to\n# just their first (leading) number.\n#\n# The document is a single-column table where every row holds exactly\n# one data point in its lone cell, so each edit is addressed directly\n# by (row, column) -- setting Cell.Range.Text replaces the whole cell\n# paragraph's text (collapsing any extra tab-separated runs) while\n# keeping the existing run formatting (font/size) intact.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"881\"\n$t.Cell(7, 1).Range.Text = \"0.05588\"\n$t.Cell(8, 1).Range.Text = \"0.03023\"\n$t.Cell(10, 1).Range.Text = \"0.11146\"\n$t.Cell(12, 1).Range.Text = \"59.86878\"\n$t.Cell(44, 1).Range.Text = \"75.98\"\n$t.Cell(45, 1).Range.Text = \"59.87\"\n$t.Cell(46, 1).Range.Text = \"249\"\n"}
